$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (15.42578125 -> 16.42578125)
$ws.Columns.Item(1).ColumnWidth = 16.42578125

# Update cell values in column A (rows 1-33)
$ws.Cells.Item(1, 1).Value = 0.20466831670470498
$ws.Cells.Item(2, 1).Value = -0.0059999999659865466
$ws.Cells.Item(3, 1).Value = -0.0039999999682702736
$ws.Cells.Item(4, 1).Value = -0.0079999999421893619
$ws.Cells.Item(5, 1).Value = -0.0029999999650502929
$ws.Cells.Item(6, 1).Value = 0.054297471790942353
$ws.Cells.Item(7, 1).Value = -0.0099999999154127828
$ws.Cells.Item(8, 1).Value = -0.019127478994534375
$ws.Cells.Item(9, 1).Value = -0.001999999959275911
$ws.Cells.Item(10, 1).Value = -0.0019999999604092267
$ws.Cells.Item(11, 1).Value = 0.023666237308498594
$ws.Cells.Item(12, 1).Value = -0.0034999999521154734
$ws.Cells.Item(13, 1).Value = -0.0034999999514582214
$ws.Cells.Item(14, 1).Value = 0.0081357937330244212
$ws.Cells.Item(15, 1).Value = -0.00099999996442168282
$ws.Cells.Item(16, 1).Value = -0.0019999999584023875
$ws.Cells.Item(17, 1).Value = -0.0019999999575750493
$ws.Cells.Item(18, 1).Value = -0.0039999999464530589
$ws.Cells.Item(19, 1).Value = -0.0039999999767279526
$ws.Cells.Item(20, 1).Value = -0.0039999999750595094
$ws.Cells.Item(21, 1).Value = -0.0039999999747806214
$ws.Cells.Item(22, 1).Value = -0.0039999999745612413
$ws.Cells.Item(23, 1).Value = -0.0049999999609946499
$ws.Cells.Item(24, 1).Value = -0.019999999876189278
$ws.Cells.Item(25, 1).Value = -0.019999999874558583
$ws.Cells.Item(26, 1).Value = -0.059938407222960066
$ws.Cells.Item(27, 1).Value = -0.0024999999512576032
$ws.Cells.Item(28, 1).Value = -0.0019999999439122007
$ws.Cells.Item(29, 1).Value = -0.0069999999090510912
$ws.Cells.Item(30, 1).Value = -0.059999999616486388
$ws.Cells.Item(31, 1).Value = -0.0069999999025522897
$ws.Cells.Item(32, 1).Value = -0.0099999998858066874
$ws.Cells.Item(33, 1).Value = -0.0039999999183280011
